$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the PAYOUT% column (I) as Text first so assigning strings like
# "16%" are preserved as literal text instead of being auto-converted
# to a numeric percentage value.
$ws.Range("I2:I13").NumberFormat = "@"

# Data for rows 2-13 (columns A-J)
# Row 2: PL Sal / V1
$ws.Cells.Item(2,1).Value = "PL Sal"
$ws.Cells.Item(2,2).Value = "V1"
$ws.Cells.Item(2,3).Value = 9170193.81
$ws.Cells.Item(2,4).Value = 30
$ws.Cells.Item(2,5).Value = 27
$ws.Cells.Item(2,6).Value = 3
$ws.Cells.Item(2,7).Value = 161172
$ws.Cells.Item(2,8).Value = 1.76
$ws.Cells.Item(2,9).Value = "16%"
$ws.Cells.Item(2,10).Value = 25787.52

# Row 3: PL Sal / V2
$ws.Cells.Item(3,1).Value = "PL Sal"
$ws.Cells.Item(3,2).Value = "V2"
$ws.Cells.Item(3,3).Value = 17749604.12
$ws.Cells.Item(3,4).Value = 52
$ws.Cells.Item(3,5).Value = 47
$ws.Cells.Item(3,6).Value = 5
$ws.Cells.Item(3,7).Value = 181585
$ws.Cells.Item(3,8).Value = 1.02
$ws.Cells.Item(3,9).Value = "17.5%"
$ws.Cells.Item(3,10).Value = 31777.375

# Row 4: PL Sal / V3
$ws.Cells.Item(4,1).Value = "PL Sal"
$ws.Cells.Item(4,2).Value = "V3"
$ws.Cells.Item(4,3).Value = 1448782.22
$ws.Cells.Item(4,4).Value = 6
$ws.Cells.Item(4,5).Value = 5
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 5000
$ws.Cells.Item(4,8).Value = 0.35
$ws.Cells.Item(4,9).Value = "20%"
$ws.Cells.Item(4,10).Value = 1000

# Row 5: PL Sal / V4
$ws.Cells.Item(5,1).Value = "PL Sal"
$ws.Cells.Item(5,2).Value = "V4"
$ws.Cells.Item(5,3).Value = 10197523.39
$ws.Cells.Item(5,4).Value = 31
$ws.Cells.Item(5,5).Value = 29
$ws.Cells.Item(5,6).Value = 2
$ws.Cells.Item(5,7).Value = 30000
$ws.Cells.Item(5,8).Value = 0.29
$ws.Cells.Item(5,9).Value = "22.5%"
$ws.Cells.Item(5,10).Value = 6750

# Row 6: PL Sal / V5
$ws.Cells.Item(6,1).Value = "PL Sal"
$ws.Cells.Item(6,2).Value = "V5"
$ws.Cells.Item(6,3).Value = 13244323.61
$ws.Cells.Item(6,4).Value = 41
$ws.Cells.Item(6,5).Value = 40
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 15000
$ws.Cells.Item(6,8).Value = 0.11
$ws.Cells.Item(6,9).Value = "25%"
$ws.Cells.Item(6,10).Value = 3750

# Row 7: PL Sal / V6
$ws.Cells.Item(7,1).Value = "PL Sal"
$ws.Cells.Item(7,2).Value = "V6"
$ws.Cells.Item(7,3).Value = 1828497.2
$ws.Cells.Item(7,4).Value = 11
$ws.Cells.Item(7,5).Value = 10
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 25000
$ws.Cells.Item(7,8).Value = 1.37
$ws.Cells.Item(7,9).Value = "32.5%"
$ws.Cells.Item(7,10).Value = 8125

# Row 8: PL Self / V1
$ws.Cells.Item(8,1).Value = "PL Self"
$ws.Cells.Item(8,2).Value = "V1"
$ws.Cells.Item(8,3).Value = 2805254.55
$ws.Cells.Item(8,4).Value = 6
$ws.Cells.Item(8,5).Value = 5
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 15000
$ws.Cells.Item(8,8).Value = 0.53
$ws.Cells.Item(8,9).Value = "16%"
$ws.Cells.Item(8,10).Value = 2400

# Row 9: PL Self / V2
$ws.Cells.Item(9,1).Value = "PL Self"
$ws.Cells.Item(9,2).Value = "V2"
$ws.Cells.Item(9,3).Value = 7479921.47
$ws.Cells.Item(9,4).Value = 15
$ws.Cells.Item(9,5).Value = 13
$ws.Cells.Item(9,6).Value = 2
$ws.Cells.Item(9,7).Value = 30000
$ws.Cells.Item(9,8).Value = 0.4
$ws.Cells.Item(9,9).Value = "17.5%"
$ws.Cells.Item(9,10).Value = 5250

# Row 10: PL Self / V3 (was V4)
$ws.Cells.Item(10,1).Value = "PL Self"
$ws.Cells.Item(10,2).Value = "V3"
$ws.Cells.Item(10,3).Value = 550961.31
$ws.Cells.Item(10,4).Value = 2
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 0
$ws.Cells.Item(10,8).Value = 0
$ws.Cells.Item(10,9).Value = "20%"
$ws.Cells.Item(10,10).Value = 0

# Row 11: PL Self / V4 (was V5)
$ws.Cells.Item(11,1).Value = "PL Self"
$ws.Cells.Item(11,2).Value = "V4"
$ws.Cells.Item(11,3).Value = 4588615.25
$ws.Cells.Item(11,4).Value = 9
$ws.Cells.Item(11,5).Value = 8
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 10000
$ws.Cells.Item(11,8).Value = 0.22
$ws.Cells.Item(11,9).Value = "22.5%"
$ws.Cells.Item(11,10).Value = 2250

# Row 12: PL Self / V5 (was V6)
$ws.Cells.Item(12,1).Value = "PL Self"
$ws.Cells.Item(12,2).Value = "V5"
$ws.Cells.Item(12,3).Value = 4249598.26
$ws.Cells.Item(12,4).Value = 11
$ws.Cells.Item(12,5).Value = 11
$ws.Cells.Item(12,6).Value = 0
$ws.Cells.Item(12,7).Value = 0
$ws.Cells.Item(12,8).Value = 0
$ws.Cells.Item(12,9).Value = "25%"
$ws.Cells.Item(12,10).Value = 0

# Row 13: PL Self / V6 (new row)
$ws.Cells.Item(13,1).Value = "PL Self"
$ws.Cells.Item(13,2).Value = "V6"
$ws.Cells.Item(13,3).Value = 1818015.45
$ws.Cells.Item(13,4).Value = 7
$ws.Cells.Item(13,5).Value = 6
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 32000
$ws.Cells.Item(13,8).Value = 1.76
$ws.Cells.Item(13,9).Value = "35%"
$ws.Cells.Item(13,10).Value = 11200

# Reset style on the I column (PAYOUT%) so the percent-like text values
# stay stored as plain text rather than being auto-converted to numeric
# percentages, and without leaving a residual numeric style behind.
$ws.Range("I2:I13").Style = "Normal"
